# Sprint backlog: add "Sprint 8 (M8)" sheet and update a few selections,
# matching the "Updated sprint backlog to include M8" commit.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Update the view/selection on a few of the existing sprint sheets
#    (these are incidental cursor-position changes left over from the
#    editing session).
# ------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Sprint 5 (M5)")
$ws5.Range("B5").Select() | Out-Null

$ws6 = $wb.Worksheets.Item("Sprint 6 (M6)")
$ws6.Range("A4").Select() | Out-Null

$ws7 = $wb.Worksheets.Item("Sprint 7 (M7)")
$ws7.Range("A1:G5").Select() | Out-Null

# ------------------------------------------------------------------
# 2. Add the new "Sprint 8 (M8)" sheet at the end of the tab strip.
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add($null, $lastSheet)
$ws8.Name = "Sprint 8 (M8)"

# Column widths matching the other sprint sheets' layout.
$ws8.Columns.Item(1).ColumnWidth = 71.99869791666667
$ws8.Columns.Item(2).ColumnWidth = 19.830729166666668

# Header row (bold + wrap, like the other sprint sheets).
$ws8.Range("A1").Value = "Tasks"
$ws8.Range("B1").Value = "Responsible"
$ws8.Range("C1").Value = "Status"
$ws8.Range("D1").Value = 1
$ws8.Range("E1").Value = 2
$ws8.Range("F1").Value = 3
$ws8.Range("G1").Value = 4
$ws8.Range("A1:G1").Font.Bold = $true
$ws8.Range("A1:G1").WrapText = $true

# Task rows.
$ws8.Range("A2").Value = "Individual sequence diagram"
$ws8.Range("B2").Value = "done"

$ws8.Range("A3").Value = "Make a ship class with given attributes (similar to resource.java) and shipyard class"

$ws8.Range("B4").Value = "Sthephen"
$ws8.Range("A4").Value = "Make a shipyard view and a controller and implement controller with GUI"

$ws8.Range("B3").Value = "Bhavesh, Naman, Pranil"

$ws8.Range("A5").Value = "Code Critique and Java Doc"
$ws8.Range("B5").Value = "Hunter"

# Final cursor position on the new sheet, and make it the active tab.
$ws8.Range("A11").Select() | Out-Null
